$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps the same shared string text, but the diff shows a new (duplicate)
# shared-string entry was introduced upstream for A1, shifting its index from 31 to 32.
$ws.Range("A1").Value = "HK_R_acc_LT"

# Updated numeric values for A2:A49 (new HK genes/reactions recomputation)
$ws.Range("A2").Value = 85.08009153318078
$ws.Range("A3").Value = 85.94965675057207
$ws.Range("A4").Value = 86.13272311212815
$ws.Range("A5").Value = 89.29061784897026
$ws.Range("A6").Value = 89.83981693363845
$ws.Range("A7").Value = 89.7025171624714
$ws.Range("A8").Value = 80.54919908466819
$ws.Range("A9").Value = 81.96796338672769
$ws.Range("A10").Value = 81.32723112128146
$ws.Range("A11").Value = 80.41189931350115
$ws.Range("A12").Value = 79.35926773455377
$ws.Range("A13").Value = 82.0137299771167
$ws.Range("A14").Value = 80.32036613272311
$ws.Range("A15").Value = 80.50343249427918
$ws.Range("A16").Value = 82.6086956521739
$ws.Range("A17").Value = 80.09153318077803
$ws.Range("A18").Value = 81.64759725400458
$ws.Range("A19").Value = 87.09382151029749
$ws.Range("A20").Value = 90.52631578947368
$ws.Range("A21").Value = 90.93821510297482
$ws.Range("A22").Value = 90.89244851258582
$ws.Range("A23").Value = 83.7070938215103
$ws.Range("A24").Value = 87.4141876430206
$ws.Range("A25").Value = 86.86498855835241
$ws.Range("A26").Value = 83.98169336384439
$ws.Range("A27").Value = 83.52402745995423
$ws.Range("A28").Value = 84.02745995423341
$ws.Range("A29").Value = 82.33409610983982
$ws.Range("A30").Value = 81.60183066361556
$ws.Range("A31").Value = 81.37299771167048
$ws.Range("A32").Value = 90.43478260869566
$ws.Range("A33").Value = 93.4096109839817
$ws.Range("A34").Value = 93.04347826086956
$ws.Range("A35").Value = 87.64302059496568
$ws.Range("A36").Value = 91.9908466819222
$ws.Range("A37").Value = 82.92906178489703
$ws.Range("A38").Value = 88.78718535469108
$ws.Range("A39").Value = 85.26315789473684
$ws.Range("A40").Value = 84.30205949656751
$ws.Range("A41").Value = 81.37299771167048
$ws.Range("A42").Value = 81.92219679633868
$ws.Range("A43").Value = 81.46453089244852
$ws.Range("A44").Value = 81.73913043478261
$ws.Range("A45").Value = 84.75972540045767
$ws.Range("A46").Value = 86.17848970251715
$ws.Range("A47").Value = 80.59496567505721
$ws.Range("A48").Value = 82.42562929061785
$ws.Range("A49").Value = 81.4187643020595
